$d = $word.ActiveDocument

# Move to the very end of the document (after the last paragraph,
# "... Aktuell ausgerüstete Waffe wird auf dem Bildschirm angezeigt").

# --- New paragraph 1: "24.05. und 25.05. " ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("24.05. und 25.05. ")

# --- New paragraph 2: three runs describing the weapon pickups ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("-Waffenpickup erstellt um Waffen im Spiel aufzusammeln. Die Shotgun schießt jetzt mehrere Kugeln und einer Kegelform, Projektile haben eine maximale Distanz die sie sich bewegen. Die Waffenpickups verteilen bei Berührung vorher zufällig erstellte Waffen")

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter(" und werden anschließend gelöscht")

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter(".")
